$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053445042572105
$ws.Range("D2").Value = 1.057432484043855
$ws.Range("E2").Value = 1.050122343835313
$ws.Range("F2").Value = 1.067646275417137
$ws.Range("I2").Value = 1.048135397382967
$ws.Range("J2").Value = 1.058462350457963
$ws.Range("K2").Value = 1.060167469247092
$ws.Range("L2").Value = 1.05287749871197
$ws.Range("M2").Value = 1.070353579898257
$ws.Range("N2").Value = 1.059965488402942
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.054574879829328
$ws.Range("D3").Value = 1.058314890067011
$ws.Range("E3").Value = 1.051091084574066
$ws.Range("F3").Value = 1.068683561301931
$ws.Range("I3").Value = 1.048472298809585
$ws.Range("J3").Value = 1.059242198020337
$ws.Range("K3").Value = 1.060863536177314
$ws.Range("L3").Value = 1.053658219422404
$ws.Range("M3").Value = 1.071206140160495
$ws.Range("N3").Value = 1.060746443438304
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.055306065315612
$ws.Range("D4").Value = 1.058885953714579
$ws.Range("E4").Value = 1.051718299149764
$ws.Range("F4").Value = 1.069355205332178
$ws.Range("I4").Value = 1.048689185117087
$ws.Range("J4").Value = 1.059746346503067
$ws.Range("K4").Value = 1.061313387368983
$ws.Range("L4").Value = 1.054163153144751
$ws.Range("M4").Value = 1.071757645003447
$ws.Range("N4").Value = 1.061251307869706
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.05561348149886
$ws.Range("D5").Value = 1.059126049700186
$ws.Range("E5").Value = 1.05198206970695
$ws.Range("F5").Value = 1.069637672230956
$ws.Range("I5").Value = 1.048780098020834
$ws.Range("J5").Value = 1.05995817927614
$ws.Range("K5").Value = 1.061502372936561
$ws.Range("L5").Value = 1.054375368665728
$ws.Range("M5").Value = 1.071989459561041
$ws.Range("N5").Value = 1.061463441469614
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.055665099559849
$ws.Range("D6").Value = 1.059166364087472
$ws.Range("E6").Value = 1.052026363192049
$ws.Range("F6").Value = 1.069685105980283
$ws.Range("I6").Value = 1.048795347092702
$ws.Range("J6").Value = 1.059993740409583
$ws.Range("K6").Value = 1.06153409671055
$ws.Range("L6").Value = 1.054410997147855
$ws.Range("M6").Value = 1.072028380009626
$ws.Range("N6").Value = 1.061499053103946
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.055310172923522
$ws.Range("D7").Value = 1.058889161807329
$ws.Range("E7").Value = 1.051721823313131
$ws.Range("F7").Value = 1.069358979245508
$ws.Range("I7").Value = 1.048690400945364
$ws.Range("J7").Value = 1.059749177459973
$ws.Range("K7").Value = 1.061315913121611
$ws.Range("L7").Value = 1.054165989009263
$ws.Range("M7").Value = 1.071760742670648
$ws.Range("N7").Value = 1.061254142846896
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053826854941776
$ws.Range("D8").Value = 1.057730678945159
$ws.Range("E8").Value = 1.050449656682881
$ws.Range("F8").Value = 1.067996737627181
$ws.Range("I8").Value = 1.048249484917916
$ws.Range("J8").Value = 1.058725999930833
$ws.Range("K8").Value = 1.060402822315078
$ws.Range("L8").Value = 1.05314139760647
$ws.Range("M8").Value = 1.070641739423493
$ws.Range("N8").Value = 1.060229512288304
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.051213837420469
$ws.Range("D9").Value = 1.055689967251405
$ws.Range("E9").Value = 1.048210809738729
$ws.Range("F9").Value = 1.065599751873853
$ws.Range("I9").Value = 1.047464022821484
$ws.Range("J9").Value = 1.056919464727532
$ws.Range("K9").Value = 1.058789623890217
$ws.Range("L9").Value = 1.051334062129796
$ws.Range("M9").Value = 1.068668701440304
$ws.Range("N9").Value = 1.058420411597802
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0494723076917
$ws.Range("D10").Value = 1.054329962482468
$ws.Range("E10").Value = 1.046720180436202
$ws.Range("F10").Value = 1.064004091853109
$ws.Range("I10").Value = 1.046934658059417
$ws.Range("J10").Value = 1.055712702454571
$ws.Range("K10").Value = 1.057711323692818
$ws.Range("L10").Value = 1.050127909295165
$ws.Range("M10").Value = 1.067352530958221
$ws.Range("N10").Value = 1.057211935583996
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04871831036829
$ws.Range("D11").Value = 1.053741178185853
$ws.Range("E11").Value = 1.046075180379
$ws.Range("F11").Value = 1.063313706437676
$ws.Range("I11").Value = 1.046704079180425
$ws.Range("J11").Value = 1.05518958793033
$ws.Range("K11").Value = 1.057243734683619
$ws.Range("L11").Value = 1.04960533106368
$ws.Range("M11").Value = 1.066782421750929
$ws.Range("N11").Value = 1.056688078177128
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.048438255725154
$ws.Range("D12").Value = 1.053522493532375
$ws.Range("E12").Value = 1.045835666470446
$ws.Range("F12").Value = 1.063057348480374
$ws.Range("I12").Value = 1.046618227409451
$ws.Range("J12").Value = 1.054995192670445
$ws.Range("K12").Value = 1.057069949294851
$ws.Range("L12").Value = 1.049411176046962
$ws.Range("M12").Value = 1.066570627824491
$ws.Range("N12").Value = 1.056493406853677
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.048498327799847
$ws.Range("D13").Value = 1.05356940140639
$ws.Range("E13").Value = 1.045887039949243
$ws.Range("F13").Value = 1.063112334426589
$ws.Range("I13").Value = 1.046636652162076
$ws.Range("J13").Value = 1.055036895066244
$ws.Range("K13").Value = 1.057107231469511
$ws.Range("L13").Value = 1.049452825050794
$ws.Range("M13").Value = 1.06661605969952
$ws.Range("N13").Value = 1.056535168471661
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.048695160693619
$ws.Range("D14").Value = 1.053723101316689
$ws.Range("E14").Value = 1.046055380701537
$ws.Range("F14").Value = 1.063292514138218
$ws.Range("I14").Value = 1.046696986818844
$ws.Range("J14").Value = 1.055173520938903
$ws.Range("K14").Value = 1.057229371611111
$ws.Range("L14").Value = 1.049589283092959
$ws.Range("M14").Value = 1.066764915411285
$ws.Range("N14").Value = 1.056671988368731
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.048816437700288
$ws.Range("D15").Value = 1.053817803031319
$ws.Range("E15").Value = 1.046159109979797
$ws.Range("F15").Value = 1.063403539658461
$ws.Range("I15").Value = 1.046734133889296
$ws.Range("J15").Value = 1.055257689079349
$ws.Range("K15").Value = 1.057304612661127
$ws.Range("L15").Value = 1.049673353269143
$ws.Range("M15").Value = 1.066856626351298
$ws.Range("N15").Value = 1.056756276037591
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.049522349961986
$ws.Range("D16").Value = 1.05436904040718
$ws.Range("E16").Value = 1.046762996484522
$ws.Range("F16").Value = 1.064049921953218
$ws.Range("I16").Value = 1.04694993214716
$ws.Range("J16").Value = 1.055747407626985
$ws.Range("K16").Value = 1.057742341738679
$ws.Range("L16").Value = 1.050162584624607
$ws.Range("M16").Value = 1.067390363031142
$ws.Range("N16").Value = 1.057246690041735
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.049965175235382
$ws.Range("D17").Value = 1.054714845637075
$ws.Range("E17").Value = 1.047141919536994
$ws.Range("F17").Value = 1.064455526579467
$ws.Range("I17").Value = 1.047084932380129
$ws.Range("J17").Value = 1.056054440026666
$ws.Range("K17").Value = 1.058016736013689
$ws.Range("L17").Value = 1.050469384432488
$ws.Range("M17").Value = 1.067725108658994
$ws.Range("N17").Value = 1.057554158462641
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.050223477105076
$ws.Range("D18").Value = 1.054916558275414
$ws.Range("E18").Value = 1.047362982587712
$ws.Range("F18").Value = 1.064692161573459
$ws.Range("I18").Value = 1.047163544476677
$ws.Range("J18").Value = 1.056233471035766
$ws.Range("K18").Value = 1.058176720154527
$ws.Range("L18").Value = 1.050648306049893
$ws.Range("M18").Value = 1.06792034095696
$ws.Range("N18").Value = 1.057733443716309
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.050311552988034
$ws.Range("D19").Value = 1.054985338818736
$ws.Range("E19").Value = 1.047438366781741
$ws.Range("F19").Value = 1.064772856929039
$ws.Range("I19").Value = 1.047190326910506
$ws.Range("J19").Value = 1.056294506529278
$ws.Range("K19").Value = 1.058231259524644
$ws.Range("L19").Value = 1.050709308696886
$ws.Range("M19").Value = 1.067986906878844
$ws.Range("N19").Value = 1.057794565887223
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04991766331104
$ws.Range("D20").Value = 1.054677742943171
$ws.Range("E20").Value = 1.047101260165863
$ws.Range("F20").Value = 1.064412003593494
$ws.Range("I20").Value = 1.047070461709903
$ws.Range("J20").Value = 1.056021504129253
$ws.Range("K20").Value = 1.057987302870259
$ws.Range("L20").Value = 1.050436470771388
$ws.Range("M20").Value = 1.067689195604103
$ws.Range("N20").Value = 1.057521175792476
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.048637197946485
$ws.Range("D21").Value = 1.053677840080814
$ws.Range("E21").Value = 1.046005806674227
$ws.Range("F21").Value = 1.063239453443183
$ws.Range("I21").Value = 1.046679225412604
$ws.Range("J21").Value = 1.055133290435369
$ws.Range("K21").Value = 1.05719340720099
$ws.Range("L21").Value = 1.049549100877495
$ws.Range("M21").Value = 1.066721081943006
$ws.Range("N21").Value = 1.056631700733268
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047832196379305
$ws.Range("D22").Value = 1.053049254331505
$ws.Range("E22").Value = 1.045317443979272
$ws.Range("F22").Value = 1.062502698576433
$ws.Range("I22").Value = 1.0464320563429
$ws.Range("J22").Value = 1.054574330366572
$ws.Range("K22").Value = 1.056693663197931
$ws.Range("L22").Value = 1.048990908841363
$ws.Range("M22").Value = 1.066112216925753
$ws.Range("N22").Value = 1.056071946877059
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.048258935886915
$ws.Range("D23").Value = 1.053382470723523
$ws.Range("E23").Value = 1.045682320912674
$ws.Range("F23").Value = 1.062893221376061
$ws.Range("I23").Value = 1.046563197543495
$ws.Range("J23").Value = 1.054870693637366
$ws.Range("K23").Value = 1.056958642956018
$ws.Range("L23").Value = 1.049286842412853
$ws.Range("M23").Value = 1.066435004283642
$ws.Range("N23").Value = 1.05636873101769
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.049939131875216
$ws.Range("D24").Value = 1.054694508021536
$ws.Range("E24").Value = 1.047119632250508
$ws.Range("F24").Value = 1.064431669594513
$ws.Range("I24").Value = 1.047077000788361
$ws.Range("J24").Value = 1.056036386615802
$ws.Range("K24").Value = 1.058000602642765
$ws.Range("L24").Value = 1.050451343129143
$ws.Range("M24").Value = 1.067705423226812
$ws.Range("N24").Value = 1.057536079413864
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.05188927711973
$ws.Range("D25").Value = 1.056217458112367
$ws.Range("E25").Value = 1.048789264483341
$ws.Range("F25").Value = 1.066219019787183
$ws.Range("I25").Value = 1.047668092347362
$ws.Range("J25").Value = 1.057386920769382
$ws.Range("K25").Value = 1.059207173956995
$ws.Range("L25").Value = 1.051801524039535
$ws.Range("M25").Value = 1.069178922333049
$ws.Range("N25").Value = 1.058888531480849
